$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its cells as Text so that values such as
# "7.13" are not silently reinterpreted as numbers by Excel (the rest of the
# column already contains values, like "56.249.32", that cannot be parsed as
# plain numbers, so formatting the whole column as Text keeps it uniform).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.249.32"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "2.990.22"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "505.90"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "137.62"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "7.13"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "0.365"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "3.507.93"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "25.68"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "0.0000162"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "56.256.42"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "2.991.47"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").Value = "8.06"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").Value = "330.86"
$ws.Range("E21").Value = "  +2.89%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "0.494"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "64.87"
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("D25").Value = "3.118.22"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "0.0₃0937"
$ws.Range("E28").Value = "  +4.64%  "
$ws.Range("D29").Value = "6.34"
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("D30").Value = "6.89"
$ws.Range("E30").Value = "  -3.66%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "20.27"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "1.15"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "153.23"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "26.30"
$ws.Range("E37").Value = "  +7.53%  "
$ws.Range("D38").Value = "1.24"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "3.030.92"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "36.78"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "0.653"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "2.182.92"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "1.34"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").Value = "0.924"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "19.44"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  -2.30%  "
